$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E (data rows 2-51) to Text format first so that
# numeric-looking strings (prices, percentages) are preserved verbatim
# (with trailing zeros, thousand-dot separators, padding spaces, etc.)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '68.522.78'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '3.932.77'
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '483.66'
$ws.Range("E5").Value = '  +5.85%  '
$ws.Range("D6").Value = '148.91'
$ws.Range("E6").Value = '  +2.29%  '
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.725'
$ws.Range("E9").Value = '  -2.95%  '
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +8.32%  '
$ws.Range("D11").Value = '0.0000352'
$ws.Range("E11").Value = '  +11.33%  '
$ws.Range("D12").Value = '42.55'
$ws.Range("E12").Value = '  -2.97%  '
$ws.Range("D13").Value = '10.59'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("D14").Value = '4.560.25'
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = '4.002.16'
$ws.Range("D16").Value = '14.65'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '19.78'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("E19").Value = '  -2.59%  '
$ws.Range("D20").Value = '68.602.87'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '432.47'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").Value = '3.37'
$ws.Range("E22").Value = '  +3.27%  '
$ws.Range("D23").Value = '14.49'
$ws.Range("E23").Value = '  -2.59%  '
$ws.Range("D24").Value = '87.16'
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").Value = '11.22'
$ws.Range("D26").Value = '3.56'
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").Value = '38.27'
$ws.Range("E27").Value = '  +1.91%  '
$ws.Range("D28").Value = '10.40'
$ws.Range("E28").Value = '  -3.23%  '
$ws.Range("D29").Value = '5.88'
$ws.Range("E29").Value = '  +6.75%  '
$ws.Range("D30").Value = '717.87'
$ws.Range("E30").Value = '  -3.79%  '
$ws.Range("D31").Value = '13.26'
$ws.Range("E31").Value = '  -4.10%  '
$ws.Range("E32").Value = '  -5.24%  '
$ws.Range("D33").Value = '2.83'
$ws.Range("E33").Value = '  +3.12%  '
$ws.Range("D34").Value = '0.0₃0897'
$ws.Range("E34").Value = '  +31.91%  '
$ws.Range("D35").Value = '41.74'
$ws.Range("E35").Value = '  -3.96%  '
$ws.Range("D36").Value = '59.17'
$ws.Range("E36").Value = '  +2.86%  '
$ws.Range("E37").Value = '  -6.90%  '
$ws.Range("D38").Value = '5.52'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  +9.39%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0469'
$ws.Range("E41").Value = '  -1.86%  '
$ws.Range("D42").Value = '3.01'
$ws.Range("E42").Value = '  +9.53%  '
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("E44").Value = '  -4.50%  '
$ws.Range("D45").Value = '0.140'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = '146.92'
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("D51").Value = '2.85'
$ws.Range("E51").Value = '  -1.42%  '

# Restore the default (unstyled) cell style on the touched data range so
# the only observable change is the cell value/text, matching the original
# formatting (these data cells had no explicit style before the edit).
$ws.Range("D2:E51").Style = "Normal"
